# Update "last_edited_time" (column D) plus the recalculated totals for
# "properties.Đầy đủ.number" (column AC), "properties.Tổng công.number"
# (column AF), and "properties.Nửa ngày.number" (column S) on the rows whose
# underlying Notion records changed between 2024-07-24 and 2024-07-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("AC2").Value = 8
$ws.Range("AF2").Value = 8

# Row 3
$ws.Range("D3").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("AC3").Value = 25
$ws.Range("AF3").Value = 25

# Row 4
$ws.Range("D4").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("AC4").Value = 25
$ws.Range("AF4").Value = 25

# Row 5
$ws.Range("D5").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("AC5").Value = 25
$ws.Range("AF5").Value = 25

# Row 6
$ws.Range("D6").Value = "2024-07-25T15:02:00.000Z"

# Row 7
$ws.Range("D7").Value = "2024-07-25T15:02:00.000Z"

# Row 8
$ws.Range("D8").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("AC8").Value = 25
$ws.Range("AF8").Value = 25

# Row 9
$ws.Range("D9").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("AC9").Value = 25
$ws.Range("AF9").Value = 25

# Row 10
$ws.Range("D10").Value = "2024-07-25T15:02:00.000Z"

# Row 11
$ws.Range("D11").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("AC11").Value = 24
$ws.Range("AF11").Value = 25

# Row 12
$ws.Range("D12").Value = "2024-07-25T15:02:00.000Z"

# Row 13
$ws.Range("D13").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("AC13").Value = 26
$ws.Range("AF13").Value = 28.5

# Row 14
$ws.Range("D14").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("S14").Value = 1
$ws.Range("AF14").Value = 24.5

# Row 15
$ws.Range("D15").Value = "2024-07-25T15:02:00.000Z"

# Row 16
$ws.Range("D16").Value = "2024-07-25T15:02:00.000Z"

# Row 17
$ws.Range("D17").Value = "2024-07-25T15:02:00.000Z"

# Row 18
$ws.Range("D18").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("AC18").Value = 19
$ws.Range("AF18").Value = 22

# Row 19
$ws.Range("D19").Value = "2024-07-25T15:02:00.000Z"

# Row 20
$ws.Range("D20").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("AC20").Value = 25
$ws.Range("AF20").Value = 25

# Row 21
$ws.Range("D21").Value = "2024-07-25T15:02:00.000Z"

# Row 22
$ws.Range("D22").Value = "2024-07-25T15:02:00.000Z"
